# Updated cryptos list on Sun Jun 11 15:57:45 UTC 2023 with GitHub Actions
#
# The Price (D) column holds plain text in the source workbook (e.g. values
# like "25.791.38" use '.' as a thousands separator, which Excel would
# otherwise happily reinterpret as a number such as 235.45 -> 235.45 (float)
# and lose the original text representation / trailing zeros). Force the
# Price column to Text before writing so values round-trip as strings, then
# restore the cells to the default ("Normal"/General) style so no stray
# formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.779.93"

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.743.87"
$ws.Range("E3").Value = "  +0.24%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "235.45"
$ws.Range("E5").Value = "  -0.06%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.05%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.5086"
$ws.Range("E7").Value = "  +3.80%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2649"
$ws.Range("E8").Value = "  +3.74%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +1.88%  "

# Row 10 - WrappedEther
$ws.Range("D10").Value = "1.749.47"
$ws.Range("E10").Value = "  +0.37%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.06931"
$ws.Range("E11").Value = "  +1.65%  "

# Row 12 - Solana
$ws.Range("D12").Value = "15.24"
$ws.Range("E12").Value = "  +3.03%  "

# Row 13 - Polygon
$ws.Range("D13").Value = "0.6187"
$ws.Range("E13").Value = "  +8.94%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.92%  "

# Row 15 - Litecoin
$ws.Range("D15").Value = "77.29"
$ws.Range("E15").Value = "  +1.54%  "

# Row 16 - BinanceUSD
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.08%  "

# Row 17 - Dai
$ws.Range("E17").Value = "  -0.05%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.802.27"
$ws.Range("E18").Value = "  +0.36%  "

# Row 19 - Avalanche
$ws.Range("D19").Value = "11.50"
$ws.Range("E19").Value = "  +1.93%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.000006574"
$ws.Range("E20").Value = "  +0.41%  "

# Row 21 - WrappedliquidstakedEther2.0
$ws.Range("D21").Value = "1.975.56"
$ws.Range("E21").Value = "  +0.57%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.042"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23 - Cosmos
$ws.Range("D23").Value = "8.213"
$ws.Range("E23").Value = "  +3.99%  "

# Row 24 - Chainlink
$ws.Range("D24").Value = "5.114"
$ws.Range("E24").Value = "  +1.95%  "

# Row 25 - Monero
$ws.Range("D25").Value = "135.84"
$ws.Range("E25").Value = "  -0.44%  "

# Row 26 - Toncoin
$ws.Range("E26").Value = "  -0.66%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "14.98"
$ws.Range("E27").Value = "  +2.38%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").Value = "1.770"
$ws.Range("E28").Value = "  -2.46%  "

# Row 29 - BitcoinCash
$ws.Range("D29").Value = "102.19"
$ws.Range("E29").Value = "  +0.77%  "

# Row 30 - Stellar
$ws.Range("D30").Value = "0.08184"
$ws.Range("E30").Value = "  +2.92%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "3.656"
$ws.Range("E31").Value = "  -1.83%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.358"
$ws.Range("E32").Value = "  -0.33%  "

# Row 33 - Hedera
$ws.Range("D33").Value = "0.04363"
$ws.Range("E33").Value = "  -0.30%  "

# Row 34 - HuobiToken
$ws.Range("E34").Value = "  +1.30%  "

# Row 35 - ARBITRUM
$ws.Range("E35").Value = "  +1.77%  "

# Row 36 - ImmutableX
$ws.Range("D36").Value = "0.5949"
$ws.Range("E36").Value = "  +0.32%  "

# Row 37 - MXToken
$ws.Range("D37").Value = "2.613"
$ws.Range("E37").Value = "  -1.77%  "

# Row 38 - VeChain
$ws.Range("D38").Value = "0.01552"
$ws.Range("E38").Value = "  +2.91%  "

# Row 39 - PaxDollar
$ws.Range("D39").Value = "1.000"
$ws.Range("E39").Value = "  -0.06%  "

# Row 40 - RenderToken
$ws.Range("D40").Value = "1.903"
$ws.Range("E40").Value = "  +1.37%  "

# Row 41 - Quant
$ws.Range("D41").Value = "101.34"
$ws.Range("E41").Value = "  +0.13%  "

# Row 42 - was TrustWalletToken, now TheSandbox
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.3804"
$ws.Range("E42").Value = "  +2.33%  "

# Row 43 - was TheSandbox, now TrustWalletToken
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.7445"
$ws.Range("E43").Value = "  +0.29%  "

# Row 44 - FraxShare
$ws.Range("D44").Value = "4.868"
$ws.Range("E44").Value = "  -5.28%  "

# Row 45 - Cronos
$ws.Range("D45").Value = "0.05483"
$ws.Range("E45").Value = "  +5.05%  "

# Row 46 - Algorand
$ws.Range("D46").Value = "0.1091"
$ws.Range("E46").Value = "  +2.10%  "

# Row 47 - Aptos
$ws.Range("D47").Value = "5.902"
$ws.Range("E47").Value = "  +2.18%  "

# Row 48 - Elrond
$ws.Range("D48").Value = "29.93"
$ws.Range("E48").Value = "  -0.23%  "

# Row 49 - Aave
$ws.Range("D49").Value = "52.27"
$ws.Range("E49").Value = "  +0.78%  "

# Row 50 - USDD
$ws.Range("D50").Value = "1.001"
$ws.Range("E50").Value = "  +0.19%  "

# Row 51 - TrueUSD
$ws.Range("D51").Value = "0.9990"
$ws.Range("E51").Value = "  +0.02%  "

# Restore the Price column's formatting to the default so no stray
# text-format styling is left on the cells.
$priceRange.NumberFormat = "General"
$priceRange.Style = "Normal"
